$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: repurpose as the PUT /Cars endpoint (update a car)
$ws.Range("F5").Value = "Updates a car given a car in the body, based on matched rego"
$ws.Range("B5").Value = "PUT"
$ws.Range("D5").Value = $ws.Range("D4").Value2
$ws.Range("D5").WrapText = $True
$ws.Range("E5").Value = $ws.Range("E4").Value2
$ws.Range("E5").WrapText = $True
$ws.Rows.Item(5).RowHeight = 38.1

# Row 6: new DELETE /Cars endpoint (remove a car)
$ws.Range("A6").Value = $ws.Range("A5").Value2
$ws.Range("B6").Value = "DELETE"
$ws.Range("D6").Value = $ws.Range("D4").Value2
$ws.Range("D6").WrapText = $True
$ws.Range("E6").Value = $ws.Range("E4").Value2
$ws.Range("E6").WrapText = $True
$ws.Rows.Item(6).RowHeight = 41.1

# Row 3 (GET /Cars/{Rego}): document the Rego parameter
$ws.Range("C3").Value = "Rego"

# Restore original row heights on untouched rows (minor font-metric rounding
# from resave is not part of the content edit)
$ws.Rows.Item(2).RowHeight = 148
$ws.Rows.Item(4).RowHeight = 107

# Move the active selection to reflect where editing ended
$ws.Range("E6").Select()
